$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: add new columns P1, Q1 (continuing the numeric sequence),
# copying the header cell format (bold, border, centered) from O1.
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap values in columns I/K/M/O and add new columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2
    $ws.Cells.Item($r, 16).Value = 2  # P: new
    $ws.Cells.Item($r, 17).Value = 2  # Q: new
}
